# Apply the "break out stock.yaml completed" edit:
#  1. Fix up column R ("backup") so that every row where R does not yet
#     mirror column Q ("detect_structure") gets corrected to match Q.
#     (This also fixes the two rows that were empty/inline-string cells,
#     which simply become numeric 0 since Q is 0 there.)
#  2. Flip the one-off isPivot flag on row 1146 (O1146: 0 -> 2).
#  3. Append six new weekly rows (1150-1155) of OHLCV + derived data,
#     leaving their "backup" (R) column unset, matching the source rows
#     immediately preceding them before this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Backfill column R from column Q wherever they disagree ----------
$lastRow = 1149
for ($r = 2; $r -le $lastRow; $r++) {
    $qCell = $ws.Cells.Item($r, 17)   # Q = detect_structure
    $rCell = $ws.Cells.Item($r, 18)   # R = backup
    $qVal = $qCell.Value2
    $rVal = $rCell.Value2
    if ($rVal -ne $qVal) {
        $rCell.Value = $qVal
    }
}

# --- 2. One-off fix: row 1146 isPivot flag -------------------------------
$ws.Cells.Item(1146, 15).Value = 2

# --- 3. Append the six new rows (1150-1155) ------------------------------
$newRows = @(
    @(45474, 2634.75,              2844,               2619.10009765625,   2669.5,             2669.5,             3729039, 2024, 7, 1,  0, 0, 0, 27, 1, 0, 0),
    @(45481, 2670.10009765625,     2740,               2586.5,             2686.14990234375,   2686.14990234375,   2172918, 2024, 7, 8,  0, 0, 0, 28, 0, 0, 0),
    @(45488, 2690,                 2744.25,            2596.14990234375,   2619.800048828125,  2619.800048828125,  1389155, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 2621,                 2699.89990234375,   2556.800048828125,  2614.14990234375,   2614.14990234375,   2036931, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 2640,                 2691.89990234375,   2431,               2435.300048828125,  2435.300048828125,  4064642, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 2415,                 2443.39990234375,   2336,               2351.550048828125,  2351.550048828125,  1863002, 2024, 8, 5,  0, 0, 0, 32, 0, 0, 0)
)

$rowIndex = 1150
foreach ($rowVals in $newRows) {
    for ($c = 1; $c -le 17; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value = $rowVals[$c - 1]
    }
    # Column A carries the same date/time number format as the rest of
    # the sheet (style index 2, numFmtId 165 = YYYY-MM-DD HH:MM:SS).
    $ws.Cells.Item($rowIndex, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    # Column R ("backup") is intentionally left blank on these new rows.
    $rowIndex++
}
